# Insert a new row at position 44 (pushes existing rows 44..164 down to 45..165)
# and populate it with the new weekly price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = 5
$ws.Range("B44").Value = "Macroferia Regional de Talca"
$ws.Range("C44").Value = "Maule"
$ws.Range("D44").Value = 44998
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100108
$ws.Range("H44").Value = "Tropicales y subtropicales"
$ws.Range("I44").Value = 100108002
$ws.Range("J44").Value = "Mango"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 248
$ws.Range("N44").Value = 7000
$ws.Range("O44").Value = 7000
$ws.Range("P44").Value = 7000
$ws.Range("Q44").Value = "$/bandeja 4 kilos"
$ws.Range("R44").Value = "Perú"
$ws.Range("S44").Value = 1750
$ws.Range("T44").Value = 4
